$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.891504666666667
$ws.Range("H2").Value = 8.674514
$ws.Range("I2").Value = 0.1213590456377548
$ws.Range("J2").Value = 0.1213590456377548
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 26.12444933333333
$ws.Range("N2").Value = 78.37334799999999
$ws.Range("O2").Value = 0.7238861157526749
$ws.Range("P2").Value = 0.7238861157526749
$ws.Range("Q2").Value = 75.53896716143021
$ws.Range("R2").Value = 679.8507044528719
$ws.Range("S2").Value = 0.08785012815816591
$ws.Range("T2").Value = 0.08785012815816591

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.891504666666667
$ws.Range("H3").Value = 8.674514
$ws.Range("I3").Value = 0.1213590456377548
$ws.Range("J3").Value = 0.1213590456377548
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.818542
$ws.Range("N3").Value = 11.455626
$ws.Range("O3").Value = 0.1058085282850919
$ws.Range("P3").Value = 0.1058085282850919
$ws.Range("Q3").Value = 11.04133201286267
$ws.Range("R3").Value = 99.37198811576401
$ws.Range("S3").Value = 0.01284082201301414
$ws.Range("T3").Value = 0.01284082201301414

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.891504666666667
$ws.Range("H4").Value = 8.674514
$ws.Range("I4").Value = 0.1213590456377548
$ws.Range("J4").Value = 0.1213590456377548
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.146179
$ws.Range("N4").Value = 18.438537
$ws.Range("O4").Value = 0.1703053559622332
$ws.Range("P4").Value = 0.1703053559622332
$ws.Range("Q4").Value = 17.77170526066867
$ws.Range("R4").Value = 159.945347346018
$ws.Range("S4").Value = 0.02066809546657473
$ws.Range("T4").Value = 0.02066809546657473

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.04042966666667
$ws.Range("H5").Value = 36.121289
$ws.Range("I5").Value = 0.505347637947847
$ws.Range("J5").Value = 0.505347637947847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 26.12444933333333
$ws.Range("N5").Value = 78.37334799999999
$ws.Range("O5").Value = 0.7238861157526749
$ws.Range("P5").Value = 0.7238861157526749
$ws.Range("Q5").Value = 314.5495947783969
$ws.Range("R5").Value = 2830.946353005572
$ws.Range("S5").Value = 0.365814138738856
$ws.Range("T5").Value = 0.365814138738856

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.04042966666667
$ws.Range("H6").Value = 36.121289
$ws.Range("I6").Value = 0.505347637947847
$ws.Range("J6").Value = 0.505347637947847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.818542
$ws.Range("N6").Value = 11.455626
$ws.Range("O6").Value = 0.1058085282850919
$ws.Range("P6").Value = 0.1058085282850919
$ws.Range("Q6").Value = 45.97688638021268
$ws.Range("R6").Value = 413.7919774219141
$ws.Range("S6").Value = 0.05347008984360917
$ws.Range("T6").Value = 0.05347008984360917

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.04042966666667
$ws.Range("H7").Value = 36.121289
$ws.Range("I7").Value = 0.505347637947847
$ws.Range("J7").Value = 0.505347637947847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.146179
$ws.Range("N7").Value = 18.438537
$ws.Range("O7").Value = 0.1703053559622332
$ws.Range("P7").Value = 0.1703053559622332
$ws.Range("Q7").Value = 74.00263596824368
$ws.Range("R7").Value = 666.0237237141931
$ws.Range("S7").Value = 0.08606340936538184
$ws.Range("T7").Value = 0.08606340936538184

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.894099000000001
$ws.Range("H8").Value = 26.682297
$ws.Range("I8").Value = 0.3732933164143983
$ws.Range("J8").Value = 0.3732933164143982
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 26.12444933333333
$ws.Range("N8").Value = 78.37334799999999
$ws.Range("O8").Value = 0.7238861157526749
$ws.Range("P8").Value = 0.7238861157526749
$ws.Range("Q8").Value = 232.3534386911507
$ws.Range("R8").Value = 2091.180948220356
$ws.Range("S8").Value = 0.270221848855653
$ws.Range("T8").Value = 0.270221848855653

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.894099000000001
$ws.Range("H9").Value = 26.682297
$ws.Range("I9").Value = 0.3732933164143983
$ws.Range("J9").Value = 0.3732933164143982
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.818542
$ws.Range("N9").Value = 11.455626
$ws.Range("O9").Value = 0.1058085282850919
$ws.Range("P9").Value = 0.1058085282850919
$ws.Range("Q9").Value = 33.962490583658
$ws.Range("R9").Value = 305.662415252922
$ws.Range("S9").Value = 0.03949761642846863
$ws.Range("T9").Value = 0.03949761642846863

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.894099000000001
$ws.Range("H10").Value = 26.682297
$ws.Range("I10").Value = 0.3732933164143983
$ws.Range("J10").Value = 0.3732933164143982
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.146179
$ws.Range("N10").Value = 18.438537
$ws.Range("O10").Value = 0.1703053559622332
$ws.Range("P10").Value = 0.1703053559622332
$ws.Range("Q10").Value = 54.664724497721
$ws.Range("R10").Value = 491.982520479489
$ws.Range("S10").Value = 0.06357385113027667
$ws.Range("T10").Value = 0.06357385113027665
